$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '66.343.76'
$ws.Range('E2').Value = '  -1.21%  '

# Row 3
$ws.Range('D3').Value = '3.546.19'
$ws.Range('E3').Value = '  +0.74%  '

# Row 4
$ws.Range('E4').Value = '  +0.06%  '

# Row 5
$ws.Range('D5').Value = "'609.46"
$ws.Range('D5').ClearFormats() | Out-Null
$ws.Range('E5').Value = '  -0.34%  '

# Row 6
$ws.Range('D6').Value = "'144.47"
$ws.Range('D6').ClearFormats() | Out-Null
$ws.Range('E6').Value = '  -2.73%  '

# Row 7
$ws.Range('D7').Value = '3.545.20'
$ws.Range('E7').Value = '  +0.70%  '

# Row 8
$ws.Range('E8').Value = '  -0.12%  '

# Row 9
$ws.Range('E9').Value = '  +0.20%  '

# Row 10
$ws.Range('D10').Value = "'8.13"
$ws.Range('D10').ClearFormats() | Out-Null
$ws.Range('E10').Value = '  -0.02%  '

# Row 11
$ws.Range('E11').Value = '  -3.97%  '

# Row 12
$ws.Range('E12').Value = '  -3.00%  '

# Row 13
$ws.Range('D13').Value = '4.142.40'
$ws.Range('E13').Value = '  +0.66%  '

# Row 14
$ws.Range('E14').Value = '  -4.60%  '

# Row 15
$ws.Range('D15').Value = "'30.18"
$ws.Range('D15').ClearFormats() | Out-Null
$ws.Range('E15').Value = '  -5.41%  '

# Row 16
$ws.Range('D16').Value = '3.541.31'
$ws.Range('E16').Value = '  +0.61%  '

# Row 17
$ws.Range('D17').Value = '66.414.59'
$ws.Range('E17').Value = '  -1.12%  '

# Row 18
$ws.Range('E18').Value = '  -0.85%  '

# Row 19
$ws.Range('D19').Value = "'10.95"
$ws.Range('D19').ClearFormats() | Out-Null
$ws.Range('E19').Value = '  +1.20%  '

# Row 20
$ws.Range('E20').Value = '  -2.94%  '

# Row 21
$ws.Range('D21').Value = "'14.93"
$ws.Range('D21').ClearFormats() | Out-Null
$ws.Range('E21').Value = '  -3.21%  '

# Row 22
$ws.Range('D22').Value = "'427.13"
$ws.Range('D22').ClearFormats() | Out-Null
$ws.Range('E22').Value = '  -2.50%  '

# Row 23
$ws.Range('E23').Value = '  -1.33%  '

# Row 24
$ws.Range('D24').Value = "'78.99"
$ws.Range('D24').ClearFormats() | Out-Null
$ws.Range('E24').Value = '  -1.03%  '

# Row 25
$ws.Range('D25').Value = '3.680.79'
$ws.Range('E25').Value = '  +0.60%  '

# Row 26
$ws.Range('E26').Value = '  -0.03%  '

# Row 27
$ws.Range('E27').Value = '  -0.57%  '

# Row 28
$ws.Range('D28').Value = "'8.13"
$ws.Range('D28').ClearFormats() | Out-Null
$ws.Range('E28').Value = '  -1.77%  '

# Row 29
$ws.Range('D29').Value = "'9.21"
$ws.Range('D29').ClearFormats() | Out-Null
$ws.Range('E29').Value = '  -6.25%  '

# Row 30
$ws.Range('E30').Value = '  -1.86%  '

# Row 31
$ws.Range('D31').Value = "'0.998"
$ws.Range('D31').ClearFormats() | Out-Null
$ws.Range('E31').Value = '  -0.14%  '

# Row 32
$ws.Range('E32').Value = '  -7.02%  '

# Row 33
$ws.Range('E33').Value = '  -4.14%  '

# Row 34
$ws.Range('D34').Value = "'25.34"
$ws.Range('D34').ClearFormats() | Out-Null
$ws.Range('E34').Value = '  -1.08%  '

# Row 35
$ws.Range('D35').Value = '3.533.05'
$ws.Range('E35').Value = '  +0.61%  '

# Row 36
$ws.Range('E36').Value = '  -0.02%  '

# Row 37
$ws.Range('D37').Value = "'1.76"
$ws.Range('D37').ClearFormats() | Out-Null
$ws.Range('E37').Value = '  -3.09%  '

# Row 38
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = "'5.65"
$ws.Range('D38').ClearFormats() | Out-Null
$ws.Range('E38').Value = '  -5.75%  '

# Row 39
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').Value = "'7.82"
$ws.Range('D39').ClearFormats() | Out-Null
$ws.Range('E39').Value = '  -3.11%  '

# Row 40
$ws.Range('D40').Value = "'0.999"
$ws.Range('D40').ClearFormats() | Out-Null
$ws.Range('E40').Value = '  +0.09%  '

# Row 41
$ws.Range('D41').Value = "'173.31"
$ws.Range('D41').ClearFormats() | Out-Null
$ws.Range('E41').Value = '  -0.99%  '

# Row 42
$ws.Range('D42').Value = "'0.0859"
$ws.Range('D42').ClearFormats() | Out-Null
$ws.Range('E42').Value = '  -4.41%  '

# Row 43
$ws.Range('D43').Value = "'5.28"
$ws.Range('D43').ClearFormats() | Out-Null
$ws.Range('E43').Value = '  -2.56%  '

# Row 44
$ws.Range('D44').Value = "'0.894"
$ws.Range('D44').ClearFormats() | Out-Null
$ws.Range('E44').Value = '  -0.23%  '

# Row 45
$ws.Range('E45').Value = '  -7.67%  '

# Row 46
$ws.Range('D46').Value = "'45.60"
$ws.Range('D46').ClearFormats() | Out-Null
$ws.Range('E46').Value = '  -1.46%  '

# Row 47
$ws.Range('E47').Value = '  -2.54%  '

# Row 48
$ws.Range('D48').Value = "'26.03"
$ws.Range('D48').ClearFormats() | Out-Null
$ws.Range('E48').Value = '  -7.98%  '

# Row 49
$ws.Range('E49').Value = '  -2.51%  '

# Row 50
$ws.Range('D50').Value = "'7.12"
$ws.Range('D50').ClearFormats() | Out-Null
$ws.Range('E50').Value = '  -4.94%  '

# Row 51
$ws.Range('D51').Value = "'0.946"
$ws.Range('D51').ClearFormats() | Out-Null
$ws.Range('E51').Value = '  -5.63%  '

Write-Output "Updated cryptos list"
